# Add TENDERS ISIC with AS-39 and AS-40 submarine tenders
# Rewrites the vessel data block (rows 2-61) with updated past-due /
# extension counts and the two new AS-39 / AS-40 tender rows inserted
# ahead of the existing submarine hull numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("SSN-775", 0, 1),
  @("SSN-776", 6, 4),
  @("SSN-777", 2, 2),
  @("SSN-786", 3, 2),
  @("SSN-792", 3, 1),
  @("SSN-794", 7, 2),
  @("SSN-762", 8, 2),
  @("SSN-763", 4, 3),
  @("SSN-766", 1, 1),
  @("SSN-771", 1, 2),
  @("SSN-772", 0, 2),
  @("SSN-773", 0, 0),
  @("SSN-752", 0, 1),
  @("SSN-754", 1, 1),
  @("SSN-758", 7, 4),
  @("SSN-767", 4, 3),
  @("SSN-722", 3, 2),
  @("SSN-759", 1, 2),
  @("SSN-760", 0, 1),
  @("SSN-761", 0, 2),
  @("SSN-783", 2, 0),
  @("SSBN-730", 1, 1),
  @("SSBN-731", 10, 2),
  @("SSBN-733", 1, 2),
  @("SSBN-739", 20, 6),
  @("SSBN-741", 2, 1),
  @("SSBN-743", 1, 1),
  @("SSGN-726", 1, 2),
  @("SSGN-727", 1, 2),
  @("SSBN-735", 0, 0),
  @("SSBN-737", 0, 1),
  @("SSN-21", 2, 2),
  @("SSN-22", 3, 0),
  @("SSN-23", 3, 1),
  @("AS-39", 1, 1),
  @("AS-40", 0, 2),
  @("SSN-774", 1, 2),
  @("SSN-778", 3, 0),
  @("SSN-779", 0, 0),
  @("SSN-780", 0, 1),
  @("SSN-781", 1, 2),
  @("SSN-795", 6, 4),
  @("SSN-753", 0, 1),
  @("SSN-756", 2, 2),
  @("SSN-764", 1, 2),
  @("SSN-765", 2, 1),
  @("SSN-791", 0, 0),
  @("SSN-793", 3, 2),
  @("SSN-750", 3, 0),
  @("SSN-751", 3, 0),
  @("SSN-768", 3, 2),
  @("SSN-769", 3, 0),
  @("SSN-785", 3, 0),
  @("SSGN-728", 2, 2),
  @("SSGN-729", 6, 4),
  @("SSBN-734", 6, 3),
  @("SSBN-736", 0, 0),
  @("SSBN-738", 1, 0),
  @("SSBN-740", 3, 2),
  @("SSBN-742", 2, 2)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}
